# Penalty Reward System edit (unfinished, per commit message):
# shifts each week's Week_Start_Date forward by one week and tweaks the
# MyForecast numbers on the "Forecast Comparison" sheet, then updates a
# handful of derived figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet -------------------------------------------
# Column B (Week_Start_Date) holds date-like text; prefix with an
# apostrophe so Excel keeps it as literal text instead of coercing it to a
# date serial number. Column D (MyForecast) stays numeric.

$forecastUpdates = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 7 }
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 7 }
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 7 }
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 11 }
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 11 }
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 10 }
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 9 }
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 9 }
    @{ Row = 10; Date = "2025-03-09"; Forecast = 9 }
    @{ Row = 11; Date = "2025-03-16"; Forecast = 7 }
    @{ Row = 12; Date = "2025-03-23"; Forecast = 7 }
    @{ Row = 13; Date = "2025-03-30"; Forecast = 9 }
    @{ Row = 14; Date = "2025-04-06"; Forecast = 8 }
    @{ Row = 15; Date = "2025-04-13"; Forecast = $null }
    @{ Row = 16; Date = "2025-04-20"; Forecast = 8 }
    @{ Row = 17; Date = "2025-04-27"; Forecast = $null }
)

foreach ($u in $forecastUpdates) {
    $wsForecast.Range("B$($u.Row)").Value = "'" + $u.Date
    if ($null -ne $u.Forecast) {
        $wsForecast.Range("D$($u.Row)").Value = $u.Forecast
    }
}

# --- "Summary" sheet --------------------------------------------------------

$wsSummary.Range("B2").Value  = "2022-12-25 to 2025-01-05"
$wsSummary.Range("B4").Value  = "'39"
$wsSummary.Range("B5").Value  = "'11"
$wsSummary.Range("B8").Value  = "1240 units"
$wsSummary.Range("B9").Value  = "'134"
$wsSummary.Range("B10").Value = "'70"
$wsSummary.Range("B11").Value = "'31"
$wsSummary.Range("B12").Value = "'11"
$wsSummary.Range("B13").Value = "'2025-02-02"
$wsSummary.Range("B14").Value = "'7"
$wsSummary.Range("B15").Value = "'2025-01-12"
